# Updates the cryptocurrency price (column D) and 1h volume-change (column E)
# figures for the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "68.155.83" },
    @{ Cell = "E2"; Value = "  -0.53%  " },
    @{ Cell = "D3"; Value = "2.640.13" },
    @{ Cell = "E3"; Value = "  -0.49%  " },
    @{ Cell = "E4"; Value = "  -0.05%  " },
    @{ Cell = "D5"; Value = "594.38" },
    @{ Cell = "E5"; Value = "  -0.99%  " },
    @{ Cell = "D6"; Value = "158.58" },
    @{ Cell = "E6"; Value = "  +2.25%  " },
    @{ Cell = "E7"; Value = "  +0.01%  " },
    @{ Cell = "E8"; Value = "  -1.21%  " },
    @{ Cell = "D9"; Value = "0.141" },
    @{ Cell = "E9"; Value = "  -3.05%  " },
    @{ Cell = "E10"; Value = "  -1.55%  " },
    @{ Cell = "D11"; Value = "5.25" },
    @{ Cell = "E11"; Value = "  -0.46%  " },
    @{ Cell = "E12"; Value = "  -1.92%  " },
    @{ Cell = "E13"; Value = "  -1.82%  " },
    @{ Cell = "D14"; Value = "3.119.04" },
    @{ Cell = "E14"; Value = "  -0.37%  " },
    @{ Cell = "D15"; Value = "0.0000186" },
    @{ Cell = "E15"; Value = "  -3.79%  " },
    @{ Cell = "D16"; Value = "68.006.91" },
    @{ Cell = "E16"; Value = "  -0.56%  " },
    @{ Cell = "D17"; Value = "2.646.79" },
    @{ Cell = "E17"; Value = "  -0.69%  " },
    @{ Cell = "D18"; Value = "11.32" },
    @{ Cell = "E18"; Value = "  -1.88%  " },
    @{ Cell = "D19"; Value = "359.37" },
    @{ Cell = "E19"; Value = "  -2.17%  " },
    @{ Cell = "D20"; Value = "7.32" },
    @{ Cell = "E20"; Value = "  -2.94%  " },
    @{ Cell = "E21"; Value = "  -1.66%  " },
    @{ Cell = "E22"; Value = "  -3.81%  " },
    @{ Cell = "D23"; Value = "2.07" },
    @{ Cell = "E23"; Value = "  -1.03%  " },
    @{ Cell = "D24"; Value = "74.58" },
    @{ Cell = "E24"; Value = "  +0.96%  " },
    @{ Cell = "E25"; Value = "  +0.00%  " },
    @{ Cell = "E26"; Value = "  -1.80%  " },
    @{ Cell = "D27"; Value = "2.773.00" },
    @{ Cell = "E27"; Value = "  -0.11%  " },
    @{ Cell = "E28"; Value = "  -5.15%  " },
    @{ Cell = "D29"; Value = "1.00" },
    @{ Cell = "D30"; Value = "558.57" },
    @{ Cell = "E30"; Value = "  -3.98%  " },
    @{ Cell = "D31"; Value = "7.96" },
    @{ Cell = "E31"; Value = "  -3.26%  " },
    @{ Cell = "E32"; Value = "  -4.62%  " },
    @{ Cell = "E33"; Value = "  -2.06%  " },
    @{ Cell = "E34"; Value = "  -0.03%  " },
    @{ Cell = "E35"; Value = "  -4.43%  " },
    @{ Cell = "D36"; Value = "1.55" },
    @{ Cell = "E36"; Value = "  -4.21%  " },
    @{ Cell = "D37"; Value = "159.62" },
    @{ Cell = "E37"; Value = "  -0.69%  " },
    @{ Cell = "E38"; Value = "  +0.61%  " },
    @{ Cell = "D39"; Value = "0.369" },
    @{ Cell = "E39"; Value = "  -1.86%  " },
    @{ Cell = "E40"; Value = "  -2.87%  " },
    @{ Cell = "D41"; Value = "5.28" },
    @{ Cell = "E41"; Value = "  -3.02%  " },
    @{ Cell = "E42"; Value = "  +0.38%  " },
    @{ Cell = "D43"; Value = "2.59" },
    @{ Cell = "E43"; Value = "  -4.78%  " },
    @{ Cell = "D44"; Value = "0.0₆0323" },
    @{ Cell = "E44"; Value = "  -3.31%  " },
    @{ Cell = "E45"; Value = "  +0.01%  " },
    @{ Cell = "D46"; Value = "156.68" },
    @{ Cell = "E46"; Value = "  -0.98%  " },
    @{ Cell = "E47"; Value = "  -2.27%  " },
    @{ Cell = "D48"; Value = "21.80" },
    @{ Cell = "E48"; Value = "  -1.34%  " },
    @{ Cell = "D49"; Value = "1.67" },
    @{ Cell = "E49"; Value = "  -3.24%  " },
    @{ Cell = "E50"; Value = "  -1.94%  " },
    @{ Cell = "E51"; Value = "  -1.41%  " }
)

foreach ($u in $updates) {
    $cellRef = $u.Cell
    $val = $u.Value
    $col = $cellRef.Substring(0, 1)
    $cell = $ws.Range($cellRef)
    if ($col -eq "D") {
        # Column D holds prices as plain text (e.g. "68.155.83", "1.00").
        # A leading apostrophe forces Excel to store the literal text instead
        # of re-parsing numeric-looking strings into numbers (which would
        # drop the thousands-dot grouping / trailing zeros). Re-apply the
        # Normal style afterwards so the quote-prefix flag set by that entry
        # mode doesn't linger on the cell's formatting.
        $cell.Value = "'" + $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}
